$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fields to expose in the client editor:
#  - 12 fields inserted right after "fuente" and before "producto"
#  - 3 fields inserted right after "producto" (fuente_base / usuario_cipre / contrasena)
$newBeforeProducto = @(
    "tipo_tramite",
    "capacidad",
    "plazo",
    "estado_civil",
    "tipo_vivienda",
    "ref1_nombre",
    "ref1_telefono",
    "ref1_parentesco",
    "ref2_nombre",
    "ref2_telefono",
    "ref2_parentesco",
    "antiguedad_cuenta"
)

$newAfterProducto = @(
    "fuente_base",
    "usuario_cipre",
    "contrasena"
)

# Locate the existing "producto" header column dynamically so the insert
# points land in the right place regardless of exact starting layout.
$productoCell = $ws.Rows(1).Find("producto")
$productoCol = $productoCell.Column

# Insert the "before producto" columns immediately in front of "producto",
# pushing "producto" (and anything after it) to the right.
for ($i = 0; $i -lt $newBeforeProducto.Count; $i++) {
    $ws.Columns($productoCol).Insert()
}

$firstNewCol = $productoCol
for ($i = 0; $i -lt $newBeforeProducto.Count; $i++) {
    $ws.Cells.Item(1, $firstNewCol + $i).Value = $newBeforeProducto[$i]
}

# "producto" has now shifted right by the number of inserted columns.
$productoCol = $productoCol + $newBeforeProducto.Count

# Insert the "after producto" columns immediately after "producto".
for ($i = 0; $i -lt $newAfterProducto.Count; $i++) {
    $ws.Columns($productoCol + 1).Insert()
}

for ($i = 0; $i -lt $newAfterProducto.Count; $i++) {
    $ws.Cells.Item(1, $productoCol + 1 + $i).Value = $newAfterProducto[$i]
}
